# Insert a new weekly record at row 8 (Dulce o Americano / Primera,
# Región de Arica y Parinacota, $/malla 70 unidades) and shift the
# existing rows 8-52 down to 9-53.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8..52 down to 9..53, creating a fresh blank row 8.
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new weekly observation.
$ws.Cells.Item(8, 1).Value2 = 11
$ws.Cells.Item(8, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(8, 3).Value2 = "Bíobío"
$ws.Cells.Item(8, 4).Value2 = 44462
$ws.Cells.Item(8, 5).Value2 = 8
$ws.Cells.Item(8, 6).Value2 = 100112024
$ws.Cells.Item(8, 7).Value2 = "Choclo"
$ws.Cells.Item(8, 8).Value2 = "Dulce o Americano"
$ws.Cells.Item(8, 9).Value2 = "Primera"
$ws.Cells.Item(8, 10).Value2 = 100
$ws.Cells.Item(8, 11).Value2 = 35000
$ws.Cells.Item(8, 12).Value2 = 36000
$ws.Cells.Item(8, 13).Value2 = 35500
$ws.Cells.Item(8, 14).Value2 = "$/malla 70 unidades"
$ws.Cells.Item(8, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value2 = 507
$ws.Cells.Item(8, 17).Value2 = 70
$ws.Cells.Item(8, 18).Value2 = "Hortaliza"

# Make sure the D column keeps the date number format used by the rest
# of the date column (copy style from the row above, which already has it).
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
